$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.350.03"
$ws.Range("D3").Value = "1.569.73"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'211.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "'0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'44.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").Value = "'0.245"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "'0.0894"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "1.585.81"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "28.350.83"
$ws.Range("D17").Value = "'0.513"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "'61.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").Value = "'227.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "0.0₃0680"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "'8.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'150.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'14.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'0.103"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "'6.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'0.0479"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("D33").Value = "'3.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "'3.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").Value = "1.378.52"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("E37").Value = "  -2.76%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").Value = "'2.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("D40").Value = "'0.0162"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0472"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.783"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").Value = "'62.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").Value = "'0.917"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.24%  "
$ws.Range("D49").Value = "1.707.01"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'85.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0513"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
